$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '94.451.70'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.37%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.429.48'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.75%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.79'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -5.83%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '643.14'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.27%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.48%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.406'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.59%  '

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.11%  '

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -6.20%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.428.69'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.78%  '

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.68%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.60'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.52%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.22'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.20%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.254.61'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.28%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.075.03'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.92%  '

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.93%  '

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -6.20%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.429.19'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.88%  '

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.82%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.53'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +6.12%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.501'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -5.59%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '498.40'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.17%  '

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -5.11%  '

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.93%  '

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -8.23%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '94.05'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.61%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.99'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.35%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.613.22'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.88%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '11.70'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.16%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.19%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.77'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +8.73%  '

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.53%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.50%  '

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.00%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '29.75'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.84%  '

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.80%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '547.20'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.40%  '

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.97%  '

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.54%  '

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.03%  '

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.33%  '

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +6.26%  '

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.22%  '

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.19%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.62'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.38%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.63'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.71%  '

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +5.40%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0410'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -4.34%  '

# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '55.08'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.44%  '

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.19'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.53%  '
